$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "Legacy Invoice Number" column (N) ---
# Copy the header style from the last existing header cell (M1) into N1, then set its text.
$ws.Range("M1").Copy($ws.Range("N1"))
$ws.Range("N1").Value2 = "Legacy Invoice Number"

# Update the email address and the number in row 2, and add the new column's data value.
$ws.Range("A2").Value2 = "webtest@yopmail.com"
$ws.Range("B2").Value2 = 2
$ws.Range("N2").Value2 = 123

# Give column N a custom width to fit the new header text.
$ws.Columns("N").ColumnWidth = 19.33

# --- Shift the "Invoice Status" legend box (rows 4-6) one column to the right ---
# Row 4: N4:P4 -> O4:Q4
$ws.Range("Q4").ClearContents()
$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("P4").ClearContents()
$ws.Range("O4").Copy($ws.Range("P4"))
$ws.Range("O4").ClearContents()
$ws.Range("N4").Copy($ws.Range("O4"))
$ws.Range("N4").Clear()

# Row 5: N5:O5 -> O5:P5
$ws.Range("P5").ClearContents()
$ws.Range("O5").Copy($ws.Range("P5"))
$ws.Range("O5").ClearContents()
$ws.Range("N5").Copy($ws.Range("O5"))
$ws.Range("N5").Clear()

# Row 6: N6:O6 -> O6:P6
$ws.Range("P6").ClearContents()
$ws.Range("O6").Copy($ws.Range("P6"))
$ws.Range("O6").ClearContents()
$ws.Range("N6").Copy($ws.Range("O6"))
$ws.Range("N6").Clear()

# Update the active cell selection.
$ws.Range("F19").Select()
